$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '79.832.38'
Set-TextCell $ws 'E2' '  +4.57%  '
Set-TextCell $ws 'D3' '3.205.22'
Set-TextCell $ws 'E3' '  +5.43%  '
Set-TextCell $ws 'D5' '206.31'
Set-TextCell $ws 'E5' '  +3.01%  '
Set-TextCell $ws 'D6' '636.46'
Set-TextCell $ws 'E6' '  +1.95%  '
Set-TextCell $ws 'E7' '  -0.02%  '
Set-TextCell $ws 'D8' '0.246'
Set-TextCell $ws 'E8' '  +19.22%  '
Set-TextCell $ws 'E9' '  +11.42%  '
Set-TextCell $ws 'D10' '3.201.96'
Set-TextCell $ws 'E10' '  +5.41%  '
Set-TextCell $ws 'D11' '0.621'
Set-TextCell $ws 'E11' '  +41.55%  '
Set-TextCell $ws 'D12' '0.0000250'
Set-TextCell $ws 'E12' '  +29.05%  '
Set-TextCell $ws 'E13' '  +3.35%  '
Set-TextCell $ws 'E14' '  +4.03%  '
Set-TextCell $ws 'D15' '3.787.78'
Set-TextCell $ws 'E15' '  +5.26%  '
Set-TextCell $ws 'D16' '32.47'
Set-TextCell $ws 'E16' '  +11.69%  '
Set-TextCell $ws 'D17' '79.662.82'
Set-TextCell $ws 'E17' '  +4.43%  '
Set-TextCell $ws 'D18' '3.200.66'
Set-TextCell $ws 'E18' '  +5.58%  '
Set-TextCell $ws 'D19' '14.66'
Set-TextCell $ws 'E19' '  +8.38%  '
Set-TextCell $ws 'D20' '9.50'
Set-TextCell $ws 'E20' '  +5.78%  '
Set-TextCell $ws 'D21' '2.96'
Set-TextCell $ws 'E21' '  +28.61%  '
Set-TextCell $ws 'D22' '435.09'
Set-TextCell $ws 'E22' '  +16.05%  '
Set-TextCell $ws 'E23' '  +20.15%  '
Set-TextCell $ws 'D24' '4.86'
Set-TextCell $ws 'E24' '  +11.22%  '
Set-TextCell $ws 'B25' 'Litecoin'
Set-TextCell $ws 'C25' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws 'D25' '77.82'
Set-TextCell $ws 'E25' '  +6.36%  '
Set-TextCell $ws 'B26' 'Aptos'
Set-TextCell $ws 'C26' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D26' '11.03'
Set-TextCell $ws 'E26' '  +12.43%  '
Set-TextCell $ws 'B27' 'Dai'
Set-TextCell $ws 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D27' '1.00'
Set-TextCell $ws 'E27' '  +0.52%  '
Set-TextCell $ws 'B28' 'PEPE'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D28' '0.0000120'
Set-TextCell $ws 'E28' '  +9.42%  '
Set-TextCell $ws 'B29' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D29' '9.30'
Set-TextCell $ws 'E29' '  +12.32%  '
Set-TextCell $ws 'B30' 'Binance-PegBSC-USD'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D30' '0.999'
Set-TextCell $ws 'E30' '  -0.17%  '
Set-TextCell $ws 'B31' 'Fetch.AI'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D31' '1.49'
Set-TextCell $ws 'E31' '  +6.15%  '
Set-TextCell $ws 'B32' 'Bittensor'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D32' '530.57'
Set-TextCell $ws 'E32' '  +7.79%  '
Set-TextCell $ws 'B33' 'PancakeSwap'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D33' '2.01'
Set-TextCell $ws 'E33' '  +3.54%  '
Set-TextCell $ws 'B34' 'Kaspa'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D34' '0.145'
Set-TextCell $ws 'E34' '  +25.76%  '
Set-TextCell $ws 'B35' 'EthereumClassic'
Set-TextCell $ws 'C35' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D35' '23.33'
Set-TextCell $ws 'E35' '  +12.98%  '
Set-TextCell $ws 'B36' 'Cronos'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 'D36' '0.123'
Set-TextCell $ws 'E36' '  +17.53%  '
Set-TextCell $ws 'B37' 'FirstDigitalUSD'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D37' '0.998'
Set-TextCell $ws 'E37' '  -0.07%  '
Set-TextCell $ws 'B38' 'PolygonEcosystemToken'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextCell $ws 'D38' '0.411'
Set-TextCell $ws 'E38' '  +7.25%  '
Set-TextCell $ws 'B39' 'Monero'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D39' '164.85'
Set-TextCell $ws 'E39' '  +1.38%  '
Set-TextCell $ws 'B40' 'WhiteBITCoin'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell $ws 'D40' '20.04'
Set-TextCell $ws 'E40' '  +0.04%  '
Set-TextCell $ws 'B41' 'Aave'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D41' '192.44'
Set-TextCell $ws 'E41' '  +1.12%  '
Set-TextCell $ws 'B42' 'USDe'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell $ws 'D42' '1.00'
Set-TextCell $ws 'E42' '  +0.02%  '
Set-TextCell $ws 'B43' 'RenderToken'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell $ws 'D43' '5.56'
Set-TextCell $ws 'E43' '  +8.42%  '
Set-TextCell $ws 'B44' 'Stacks'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D44' '1.82'
Set-TextCell $ws 'E44' '  +10.66%  '
Set-TextCell $ws 'B45' 'Mantle'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws 'D45' '0.808'
Set-TextCell $ws 'E45' '  +0.64%  '
Set-TextCell $ws 'B46' 'ImmutableX'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D46' '1.32'
Set-TextCell $ws 'E46' '  +4.64%  '
Set-TextCell $ws 'B47' 'OKB'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D47' '43.50'
Set-TextCell $ws 'E47' '  +3.47%  '
Set-TextCell $ws 'B48' 'dogwifhat'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell $ws 'D48' '2.61'
Set-TextCell $ws 'E48' '  +6.20%  '
Set-TextCell $ws 'B49' 'InjectiveProtocol'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D49' '25.82'
Set-TextCell $ws 'E49' '  +15.77%  '
Set-TextCell $ws 'B50' 'ARBITRUM'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D50' '0.639'
Set-TextCell $ws 'E50' '  +5.68%  '
Set-TextCell $ws 'B51' 'Filecoin'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D51' '4.22'
Set-TextCell $ws 'E51' '  +7.96%  '

Write-Host "Applied crypto list update"